# The edit re-shuffles the per-observation data among the existing data
# rows (rows 2-27) of the "Artfynd" sheet: every row keeps its constant
# columns (location, dates, validation status text, etc.) but the
# observation-specific columns move to a different row according to a
# fixed permutation. There are no fixed points, so every one of these
# columns must be fully snapshotted before anything is written back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27

# Columns whose values actually move between rows.
$cols = @(1, 2, 4, 5, 6, 7, 8, 9, 13, 14, 17, 18, 29)

# Snapshot every varying column up-front (single-column ranges keep the
# read far away from any of the constant, date-like text columns so we
# never risk Excel re-typing "2023-08-31" as a real date).
$snap = @{}
foreach ($c in $cols) {
    $colLetter = $ws.Cells.Item(1, $c).Address($false, $false) -replace '[0-9]+$', ''
    $rangeAddr = $colLetter + $firstRow + ":" + $colLetter + $lastRow
    $snap[$c] = $ws.Range($rangeAddr).Value2
}

# target row -> source row (pure permutation of rows 2..27, no fixed points)
$map = @{}
$map[2] = 4
$map[3] = 16
$map[4] = 24
$map[5] = 3
$map[6] = 9
$map[7] = 13
$map[8] = 10
$map[9] = 17
$map[10] = 21
$map[11] = 14
$map[12] = 8
$map[13] = 25
$map[14] = 19
$map[15] = 23
$map[16] = 22
$map[17] = 18
$map[18] = 7
$map[19] = 27
$map[20] = 11
$map[21] = 12
$map[22] = 6
$map[23] = 26
$map[24] = 5
$map[25] = 20
$map[26] = 2
$map[27] = 15

for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $map[$targetRow]
    $srcIdx = $sourceRow - $firstRow + 1
    foreach ($c in $cols) {
        $val = $snap[$c][$srcIdx, 1]
        $ws.Cells.Item($targetRow, $c).Value2 = $val
    }
}
